# Updated remaining queries for C3DC
# Replaces the old "id"-based join keys with the correct "study_id" /
# "participant_id" keys across every SQL query cell on the sheet, moves the
# active selection, and widens column C to fit the longer query text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-QueryText([string]$text) {
    $text = $text.Replace('std.id = prt."study.id"', 'std.study_id = prt."study.study_id"')
    $text = $text.Replace('prt.id = dgn."participant.id"', 'prt.participant_id = dgn."participant.participant_id"')
    $text = $text.Replace('prt.id = trt."participant.id"', 'prt.participant_id = trt."participant.participant_id"')
    $text = $text.Replace('prt.id = trr."participant.id"', 'prt.participant_id = trr."participant.participant_id"')
    $text = $text.Replace('prt.id = srv."participant.id"', 'prt.participant_id = srv."participant.participant_id"')
    $text = $text.Replace('std.id = rfs."study.id"', 'std.study_id = rfs."study.study_id"')
    return $text
}

$queryCells = @("C2", "B2", "B3", "B4", "B5", "B6", "B7")
foreach ($addr in $queryCells) {
    $cell = $ws.Range($addr)
    $cell.Value2 = Update-QueryText $cell.Value2
}

# Move the active cell selection from C7 to B2.
$ws.Range("B2").Select()

# Column C no longer auto-fits (bestFit) and is instead a fixed, wider width
# (~70.66 chars, the closest this engine's character-width grid can land to
# the authored 70.6640625).
$ws.Columns.Item(3).ColumnWidth = 69.83
